# Generate Report for Handoff
#
# The localization status report used to track 2 files (a .md + the
# .localization-config). It now tracks 4 files:
#   624b1b44-083f-4252-a5e3-6eb81e0a2814.png   (row 2)
#   97a4ae30-268f-460f-b9f3-b07941c42621.md    (row 3, replaces old .md)
#   c026a874-4b95-4288-baf3-a303a0db4565.png   (row 4, new)
#   .localization-config                        (row 5, shifted down)
#
# Rewritten on all 3 sheets (Overview + one detail sheet per locale).

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/7f675de1a6f3ea9525d80e75ba716cfe700fb934"
$zhBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d0c7ae36eb95ff4f0fd03d6566843713bb57961/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5225201122f561ab61e56df9c4d4daf56a9b1e9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$filePng1 = "624b1b44-083f-4252-a5e3-6eb81e0a2814.png"
$fileMd   = "97a4ae30-268f-460f-b9f3-b07941c42621.md"
$filePng2 = "c026a874-4b95-4288-baf3-a303a0db4565.png"
$config   = ".localization-config"

$ready  = "Ready for handoff"
$notLoc = "Not to be localized"

$handoffPng1Zh = "2048ef873049e13574ca285e7e969ce92dd5e879.png"
$handoffMdZh   = "97a4ae30-268f-460f-b9f3-b07941c42621.29f9e29b026afaecfba8dba4a8806b8073cf0234.zh-cn.xlf"
$handoffPng2Zh = "20b2374cf8553d019b4763519097cc6601337886.png"

$handoffPng1De = "2048ef873049e13574ca285e7e969ce92dd5e879.png"
$handoffMdDe   = "97a4ae30-268f-460f-b9f3-b07941c42621.29f9e29b026afaecfba8dba4a8806b8073cf0234.de-de.xlf"
$handoffPng2De = "20b2374cf8553d019b4763519097cc6601337886.png"

$dtZh = "2016-03-09 23:06:37"
$dtDe = "2016-03-09 23:06:43"
$emptyDt = "0001-01-01 00:00:00"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de, 4 data rows now (was 2).
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Drop every existing hyperlink on the sheet - they get re-created below in
# the final row order/number (cell values + styles are untouched by this).
$wsOv.Range("A1:C5").Hyperlinks.Delete()

$wsOv.Range("B2").Value = $ready
$wsOv.Range("C2").Value = $ready

$wsOv.Range("B3").Value = $ready
$wsOv.Range("C3").Value = $ready

$wsOv.Range("B4").Value = $ready
$wsOv.Range("C4").Value = $ready

$wsOv.Range("B5").Value = $notLoc
$wsOv.Range("C5").Value = $notLoc

$wsOv.Hyperlinks.Add($wsOv.Range("A2"), "$repoBase/e2e/$filePng1", "", "", $filePng1)
$wsOv.Hyperlinks.Add($wsOv.Range("A3"), "$repoBase/e2e/$fileMd", "", "", $fileMd)
$wsOv.Hyperlinks.Add($wsOv.Range("A4"), "$repoBase/e2e/$filePng2", "", "", $filePng2)
$wsOv.Hyperlinks.Add($wsOv.Range("A5"), "$repoBase/$config", "", "", $config)

# ---------------------------------------------------------------------------
# Per-locale detail sheets ("zh-cn", "de-de"): same row layout, locale-
# specific handoff file names / timestamps.
# ---------------------------------------------------------------------------
function Set-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$XlfBase,
        [string]$HandoffPng1,
        [string]$HandoffMd,
        [string]$HandoffPng2,
        [string]$Dt
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Range("A1:I5").Hyperlinks.Delete()

    # Row 2 - first .png
    $ws.Range("B2").Value = $ready
    $ws.Range("D2").Value = $Dt
    $ws.Range("D2").NumberFormat = $dateFmt
    $ws.Range("G2").Value = $emptyDt
    $ws.Range("G2").NumberFormat = $dateFmt
    $ws.Range("H2").Value = "IsDependency"
    $ws.Range("I2").Value = "e2e\$fileMd"

    # Row 3 - .md file (used to be the .localization-config row)
    $ws.Range("B3").Value = $ready
    $ws.Range("D3").Value = $Dt
    $ws.Range("D3").NumberFormat = $dateFmt
    $ws.Range("G3").Value = $emptyDt
    $ws.Range("G3").NumberFormat = $dateFmt
    $ws.Range("H3").Value = "Include"

    # Row 4 (new) - second .png
    $ws.Range("B4").Value = $ready
    $ws.Range("D4").Value = $Dt
    $ws.Range("D4").NumberFormat = $dateFmt
    $ws.Range("G4").Value = $emptyDt
    $ws.Range("G4").NumberFormat = $dateFmt
    $ws.Range("H4").Value = "IsDependency"
    $ws.Range("I4").Value = "e2e\$fileMd"

    # Row 5 (new) - .localization-config, shifted down from row 3
    $ws.Range("B5").Value = $notLoc
    $ws.Range("D5").Value = $emptyDt
    $ws.Range("D5").NumberFormat = $dateFmt
    $ws.Range("G5").Value = $emptyDt
    $ws.Range("G5").NumberFormat = $dateFmt
    $ws.Range("H5").Value = "Ignored"

    $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/$filePng1", "", "", $filePng1)
    $ws.Hyperlinks.Add($ws.Range("C2"), "$XlfBase/$HandoffPng1", "", "", $HandoffPng1)
    $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/e2e/$fileMd", "", "", $fileMd)
    $ws.Hyperlinks.Add($ws.Range("C3"), "$XlfBase/$HandoffMd", "", "", $HandoffMd)
    $ws.Hyperlinks.Add($ws.Range("A4"), "$repoBase/e2e/$filePng2", "", "", $filePng2)
    $ws.Hyperlinks.Add($ws.Range("C4"), "$XlfBase/$HandoffPng2", "", "", $HandoffPng2)
    $ws.Hyperlinks.Add($ws.Range("A5"), "$repoBase/$config", "", "", $config)
}

Set-LocaleSheet "zh-cn" $zhBase $handoffPng1Zh $handoffMdZh $handoffPng2Zh $dtZh
Set-LocaleSheet "de-de" $deBase $handoffPng1De $handoffMdDe $handoffPng2De $dtDe
